# Auto-generated Excel COM-interop script
# Applies numeric updates (and a few cell clears) to the Leve profit
# tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), matching
# a scheduled-runner data refresh of currentAveragePrice* / Leve*Profit* columns.

$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H88").Value = 996.8
$ws.Range("I88").Value = 725
$ws.Range("J88").Value = 1064.75
$ws.Range("K88").Value = 725
$ws.Range("L88").Value = 1064.75
$ws.Range("M88").Value = -319
$ws.Range("N88").Value = -1876.75
$ws.Range("H91").Value = 996.8
$ws.Range("I91").Value = 725
$ws.Range("J91").Value = 1064.75
$ws.Range("K91").Value = 725
$ws.Range("L91").Value = 1064.75
$ws.Range("M91").Value = 679
$ws.Range("N91").Value = -3872.75
$ws.Range("H98").Value = 741.75
$ws.Range("I98").Value = 676.2857
$ws.Range("J98").Value = 1200
$ws.Range("K98").Value = 676.2857
$ws.Range("L98").Value = 1200
$ws.Range("M98").Value = 821.7143
$ws.Range("N98").Value = -4196
$ws.Range("H112").Value = 2514.8333
$ws.Range("J112").Value = 1897.25
$ws.Range("L112").Value = 5691.75
$ws.Range("N112").Value = -7907.75
$ws.Range("H113").Value = 7419.7334
$ws.Range("I113").Value = 7188.154
$ws.Range("K113").Value = 7188.154
$ws.Range("M113").Value = -3934.154
$ws.Range("H122").Value = 741.75
$ws.Range("I122").Value = 676.2857
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 2028.8571
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = 421.1428999999998
$ws.Range("N122").Value = -8500
$ws.Range("H132").Value = 2275.75
$ws.Range("I132").Value = 1243.7142
$ws.Range("J132").Value = 9500
$ws.Range("K132").Value = 3731.1426
$ws.Range("L132").Value = 28500
$ws.Range("M132").Value = -1201.1426
$ws.Range("N132").Value = -33560
$ws.Range("H138").Value = 2710.0715
$ws.Range("I138").Value = 1543.3334
$ws.Range("K138").Value = 4630.0002
$ws.Range("M138").Value = 509.9997999999996
$ws.Range("H141").Value = 1622.9
$ws.Range("I141").Value = 1278.75
$ws.Range("J141").Value = 2999.5
$ws.Range("K141").Value = 3836.25
$ws.Range("L141").Value = 8998.5
$ws.Range("M141").Value = 1343.75
$ws.Range("N141").Value = -19358.5

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H4").Value = 491.8
$ws.Range("J4").Value = 571
$ws.Range("L4").Value = 571
$ws.Range("N4").Value = -803
$ws.Range("H61").Value = 1580.4706
$ws.Range("I61").Value = 1580.4706
$ws.Range("K61").Value = 1580.4706
$ws.Range("M61").Value = -1368.4706
$ws.Range("H121").Value = 79998
$ws.Range("J121").Value = 79998
$ws.Range("L121").Value = 79998
$ws.Range("N121").Value = -83492
$ws.Range("H132").Value = 1211.1428
$ws.Range("I132").Value = 1177.909
$ws.Range("K132").Value = 3533.727
$ws.Range("M132").Value = -1003.727
$ws.Range("H136").Value = 1580.4706
$ws.Range("I136").Value = 1580.4706
$ws.Range("K136").Value = 4741.4118
$ws.Range("M136").Value = -2191.4118

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H22").Value = 310
$ws.Range("I22").Value = 310
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 310
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -137
$ws.Range("H86").Value = 3248.2856
$ws.Range("J86").Value = 6700
$ws.Range("L86").Value = 6700
$ws.Range("N86").Value = -8946
$ws.Range("H89").Value = 3248.2856
$ws.Range("J89").Value = 6700
$ws.Range("L89").Value = 33500
$ws.Range("N89").Value = -44732
$ws.Range("H105").Value = 30304366
$ws.Range("I105").Value = 30304366
$ws.Range("K105").Value = 30304366
$ws.Range("M105").Value = -30302619
$ws.Range("N22").ClearContents()

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H4").Value = 1990
$ws.Range("I4").Value = 1990
$ws.Range("K4").Value = 1990
$ws.Range("M4").Value = -1878
$ws.Range("H88").Value = 5407.3335
$ws.Range("J88").Value = 7000
$ws.Range("L88").Value = 7000
$ws.Range("N88").Value = -7812
$ws.Range("H91").Value = 5407.3335
$ws.Range("J91").Value = 7000
$ws.Range("L91").Value = 7000
$ws.Range("N91").Value = -9808
$ws.Range("H99").Value = 3905.2
$ws.Range("I99").Value = 3999.5
$ws.Range("K99").Value = 3999.5
$ws.Range("M99").Value = -2501.5
$ws.Range("H105").Value = 2458.5
$ws.Range("I105").Value = 2098.5
$ws.Range("J105").Value = 3898.5
$ws.Range("K105").Value = 2098.5
$ws.Range("L105").Value = 3898.5
$ws.Range("M105").Value = -351.5
$ws.Range("N105").Value = -7392.5
$ws.Range("H126").Value = 3905.2
$ws.Range("I126").Value = 3999.5
$ws.Range("K126").Value = 11998.5
$ws.Range("M126").Value = -9528.5
$ws.Range("H132").Value = 1708.4117
$ws.Range("I132").Value = 1708.4117
$ws.Range("K132").Value = 5125.2351
$ws.Range("M132").Value = -2595.2351

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H68").Value = 357.75
$ws.Range("I68").Value = 357.75
$ws.Range("K68").Value = 1073.25
$ws.Range("M68").Value = -262.25
$ws.Range("H71").Value = 357.75
$ws.Range("I71").Value = 357.75
$ws.Range("K71").Value = 3219.75
$ws.Range("M71").Value = 836.25
$ws.Range("H122").Value = 699
$ws.Range("I122").Value = 348.5
$ws.Range("J122").Value = 1224.75
$ws.Range("K122").Value = 3136.5
$ws.Range("L122").Value = 11022.75
$ws.Range("M122").Value = -686.5
$ws.Range("N122").Value = -15922.75
$ws.Range("H129").Value = 2633.5557
$ws.Range("I129").Value = 1632.75
$ws.Range("J129").Value = 3434.2
$ws.Range("K129").Value = 4898.25
$ws.Range("L129").Value = 10302.6
$ws.Range("M129").Value = 101.75
$ws.Range("N129").Value = -20302.6

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H70").Value = 8999.5
$ws.Range("J70").Value = 8000
$ws.Range("L70").Value = 8000
$ws.Range("N70").Value = -8540
$ws.Range("H73").Value = 8999.5
$ws.Range("J73").Value = 8000
$ws.Range("L73").Value = 8000
$ws.Range("N73").Value = -9872
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("H126").Value = 2498.5
$ws.Range("I126").Value = 2498.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7495.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5025.5
$ws.Range("H132").Value = 84098
$ws.Range("I132").Value = 84098
$ws.Range("K132").Value = 252294
$ws.Range("M132").Value = -249764
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H22").Value = 931.25
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("H27").Value = 931.25
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("H100").Value = 8999.75
$ws.Range("I100").Value = 7332.6665
$ws.Range("K100").Value = 7332.6665
$ws.Range("M100").Value = -6791.6665
$ws.Range("H132").Value = 4671.4287
$ws.Range("I132").Value = 4671.4287
$ws.Range("K132").Value = 14014.2861
$ws.Range("M132").Value = -11484.2861
$ws.Range("M22").ClearContents()
$ws.Range("M27").ClearContents()

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("H54").Value = 9027.286
$ws.Range("I54").Value = 4800
$ws.Range("K54").Value = 4800
$ws.Range("M54").Value = -4280
$ws.Range("H115").Value = 60000
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("H132").Value = 2598.25
$ws.Range("I132").Value = 2598.25
$ws.Range("K132").Value = 7794.75
$ws.Range("M132").Value = -5264.75
$ws.Range("H136").Value = 3076.1333
$ws.Range("I136").Value = 2132
$ws.Range("K136").Value = 6396
$ws.Range("M136").Value = -3846
$ws.Range("N49").ClearContents()
$ws.Range("M115").ClearContents()
